# Commit: "adding version number to the ppt"
#
# Slide 1 ("Dream rich" / "Park Finder" title slide): the subtitle
# placeholder's text "Park Finder" gets split into two runs ("Park " +
# "Finder") and a new paragraph "v1" is appended below it.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$shape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "Subtitle 2") {
        $shape = $candidate
    }
}
if ($shape -eq $null) {
    $shape = $s.Shapes.Item(2)
}

$tr = $shape.TextFrame.TextRange
$tr.Text = "Park "
[void]$tr.InsertAfter("Finder")

$full = $shape.TextFrame.TextRange
[void]$full.InsertAfter("`rv1")
